# Regenerate orders with updated distance/size codes.
# Distance codes: D64 -> D69, D51 -> D55, D80 -> D86
# Size code:      S30 -> S31
# These substrings appear (possibly several times) inside string values in
# columns B (Condition), D (Filename_Left), E (Filename_Right), H (Distance)
# and J (Size). Every other column/value is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2

        if ($val -is [string]) {
            $newVal = $val.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
